# Update scraped "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages build and the one generated at commit 456a3b4.
#
# Sheet "展览" (worksheet 1) and sheet "全部类型" (worksheet 4) both list the
# same events, so each event's updated count needs to be written to both
# sheets (the two sheets drifted slightly apart on row F12/F7 because the
# counts were captured a few seconds apart during the site scrape).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 12
$wsExhibit.Range("F7").Value = 10706
$wsExhibit.Range("F30").Value = 1185
$wsExhibit.Range("F33").Value = 1416
$wsExhibit.Range("F37").Value = 22
$wsExhibit.Range("F38").Value = 130

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F11").Value = 12
$wsAll.Range("F12").Value = 10707
$wsAll.Range("F28").Value = 1185
$wsAll.Range("F34").Value = 1416
$wsAll.Range("F38").Value = 130
